$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Column D ("Price") values are plain text in the source data (e.g. "314.99",
# "45.667.12"). Excel auto-converts strings that look numeric into real numbers when
# assigned via .Value, which would corrupt formatting (e.g. "45.667.12" or trailing
# zeros like "1.00"/"3.00"). Prefixing with a leading apostrophe forces Excel to keep
# the value as literal text, matching the original inlineStr cell content exactly.

$ws.Range('D2').Value = '''45.667.12'
$ws.Range('E2').Value = '  +7.05%  '

$ws.Range('D3').Value = '''2.383.64'
$ws.Range('E3').Value = '  +3.40%  '

$ws.Range('E4').Value = '  +0.36%  '

$ws.Range('D5').Value = '''111.72'
$ws.Range('E5').Value = '  +6.06%  '

$ws.Range('D6').Value = '''314.99'
$ws.Range('E6').Value = '  +1.77%  '

$ws.Range('D7').Value = '''0.630'
$ws.Range('E7').Value = '  +0.66%  '

$ws.Range('E8').Value = '  -0.07%  '

$ws.Range('D9').Value = '''0.619'
$ws.Range('E9').Value = '  +2.21%  '

$ws.Range('D10').Value = '''41.25'
$ws.Range('E10').Value = '  +3.58%  '

$ws.Range('D11').Value = '''0.0923'
$ws.Range('E11').Value = '  +1.64%  '

$ws.Range('D12').Value = '''8.58'
$ws.Range('E12').Value = '  +3.44%  '

$ws.Range('D13').Value = '''0.109'
$ws.Range('E13').Value = '  +1.90%  '

$ws.Range('D14').Value = '''0.988'
$ws.Range('E14').Value = '  -0.01%  '

$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '''2.758.42'
$ws.Range('E15').Value = '  +4.10%  '

$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = '''15.62'
$ws.Range('E16').Value = '  +2.03%  '

$ws.Range('D17').Value = '''2.396.49'
$ws.Range('E17').Value = '  +4.53%  '

$ws.Range('D18').Value = '''45.637.28'
$ws.Range('E18').Value = '  +6.59%  '

$ws.Range('D19').Value = '''7.38'
$ws.Range('E19').Value = '  +0.70%  '

$ws.Range('D20').Value = '''0.0000107'
$ws.Range('E20').Value = '  +1.91%  '

$ws.Range('D21').Value = '''13.15'
$ws.Range('E21').Value = '  -4.30%  '

$ws.Range('D22').Value = '''73.93'
$ws.Range('E22').Value = '  +0.64%  '

$ws.Range('D23').Value = '''3.50'
$ws.Range('E23').Value = '  +1.25%  '

$ws.Range('D24').Value = '''262.49'
$ws.Range('E24').Value = '  -2.18%  '

$ws.Range('D25').Value = '''2.31'
$ws.Range('E25').Value = '  +3.04%  '

$ws.Range('E26').Value = '  -0.54%  '

$ws.Range('D27').Value = '''7.54'
$ws.Range('E27').Value = '  -0.91%  '

$ws.Range('D28').Value = '''11.17'
$ws.Range('E28').Value = '  +1.98%  '

$ws.Range('D29').Value = '''2.36'
$ws.Range('E29').Value = '  +2.43%  '

$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = '''22.58'
$ws.Range('E30').Value = '  +1.54%  '

$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = '''38.30'
$ws.Range('E31').Value = '  +1.07%  '

$ws.Range('D32').Value = '''0.0964'
$ws.Range('E32').Value = '  +11.24%  '

$ws.Range('D33').Value = '''171.01'
$ws.Range('E33').Value = '  +3.26%  '

$ws.Range('E34').Value = '  +4.25%  '

$ws.Range('D35').Value = '''0.131'
$ws.Range('E35').Value = '  +0.44%  '

$ws.Range('E36').Value = '  +4.43%  '

$ws.Range('D37').Value = '''4.81'
$ws.Range('E37').Value = '  +4.03%  '

$ws.Range('D38').Value = '''4.01'
$ws.Range('E38').Value = '  +10.83%  '

$ws.Range('D39').Value = '''3.00'
$ws.Range('E39').Value = '  +7.26%  '

$ws.Range('D40').Value = '''0.0357'
$ws.Range('E40').Value = '  +0.37%  '

$ws.Range('D41').Value = '''1.73'
$ws.Range('E41').Value = '  +10.36%  '

$ws.Range('D42').Value = '''101.87'
$ws.Range('E42').Value = '  -5.96%  '

$ws.Range('D43').Value = '''0.237'
$ws.Range('E43').Value = '  +3.57%  '

$ws.Range('D44').Value = '''13.31'
$ws.Range('E44').Value = '  +8.56%  '

$ws.Range('D45').Value = '''70.78'
$ws.Range('E45').Value = '  -1.02%  '

$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').Value = '''1.00'
$ws.Range('E46').Value = '  -0.21%  '

$ws.Range('B47').Value = 'ordi'
$ws.Range('C47').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D47').Value = '''85.39'
$ws.Range('E47').Value = '  +12.68%  '

$ws.Range('D48').Value = '''113.95'
$ws.Range('E48').Value = '  +1.91%  '

$ws.Range('D49').Value = '''9.43'
$ws.Range('E49').Value = '  +6.35%  '

$ws.Range('D50').Value = '''5.57'
$ws.Range('E50').Value = '  +7.57%  '

$ws.Range('D51').Value = '''1.651.57'
$ws.Range('E51').Value = '  -2.91%  '
